$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 418
$ws1.Range("F5").Value = 1246
$ws1.Range("F7").Value = 7511
$ws1.Range("F8").Value = 91
$ws1.Range("F11").Value = 8177
$ws1.Range("F14").Value = 5584
$ws1.Range("F16").Value = 2546
$ws1.Range("F17").Value = 1092
$ws1.Range("F18").Value = 4574
$ws1.Range("F19").Value = 320
$ws1.Range("F22").Value = 25
$ws1.Range("F23").Value = 458
$ws1.Range("F24").Value = 1997
$ws1.Range("F26").Value = 2715
$ws1.Range("F28").Value = 310
$ws1.Range("F29").Value = 108
$ws1.Range("F30").Value = 256
$ws1.Range("F31").Value = 626
$ws1.Range("F32").Value = 10
$ws1.Range("F33").Value = 532
$ws1.Range("F34").Value = 1607
$ws1.Range("F37").Value = 2554

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 32

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 418
$ws4.Range("F7").Value = 1246
$ws4.Range("F9").Value = 7511
$ws4.Range("F10").Value = 91
$ws4.Range("F13").Value = 8177
$ws4.Range("F16").Value = 5584
$ws4.Range("F18").Value = 2546
$ws4.Range("F19").Value = 1092
$ws4.Range("F20").Value = 4574
$ws4.Range("F21").Value = 320
$ws4.Range("F25").Value = 25
$ws4.Range("F27").Value = 458
$ws4.Range("F28").Value = 1997
$ws4.Range("F30").Value = 2715
$ws4.Range("F32").Value = 310
$ws4.Range("F33").Value = 108
$ws4.Range("F34").Value = 256
$ws4.Range("F35").Value = 32
$ws4.Range("F36").Value = 626
$ws4.Range("F37").Value = 10
$ws4.Range("F38").Value = 532
$ws4.Range("F40").Value = 1607
$ws4.Range("F43").Value = 2554
